$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "general": update summary stats (new objective value / breakdown)
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("B3").Value = 98.26098846897688
$wsGeneral.Range("B6").Value = 38.17098846897686
$wsGeneral.Range("B7").Value = 0
$wsGeneral.Range("B8").Value = 0
$wsGeneral.Range("B9").Value = 0
$wsGeneral.Range("B10").Value = 60.09

# ---------------------------------------------------------------------------
# Sheet "x": permutation values change
# ---------------------------------------------------------------------------
$wsX = $wb.Worksheets.Item("x")
$wsX.Range("B2").Value = 1
$wsX.Range("B4").Value = 6
$wsX.Range("B6").Value = 5
$wsX.Range("B7").Value = 4
$wsX.Range("B8").Value = 8
$wsX.Range("B9").Value = 12
$wsX.Range("B10").Value = 2
$wsX.Range("B11").Value = 10
$wsX.Range("B12").Value = 9
$wsX.Range("B13").Value = 13
$wsX.Range("B14").Value = 11

# ---------------------------------------------------------------------------
# Sheet "U"
# ---------------------------------------------------------------------------
$wsU = $wb.Worksheets.Item("U")
$wsU.Range("B4").Value = 2
$wsU.Range("B6").Value = 3
$wsU.Range("B10").Value = 3

# ---------------------------------------------------------------------------
# Sheet "TBar"
# ---------------------------------------------------------------------------
$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Range("B3").Value = 20
$wsTBar.Range("B4").Value = 24.34885042281291
$wsTBar.Range("B5").Value = 10
$wsTBar.Range("B6").Value = 20
$wsTBar.Range("B7").Value = 27.06506101847739
$wsTBar.Range("B8").Value = 20
$wsTBar.Range("B9").Value = 20.60033324079215
$wsTBar.Range("B10").Value = 22.31224998648502
$wsTBar.Range("B11").Value = 24.76592070603971
$wsTBar.Range("B12").Value = 22.61192465059683
$wsTBar.Range("B13").Value = 26.71671453559703
$wsTBar.Range("B14").Value = 28.25017704655228
$wsTBar.Range("B15").Value = 27.87444125446785

# ---------------------------------------------------------------------------
# Sheet "y": only the header row remains
# ---------------------------------------------------------------------------
$wsY = $wb.Worksheets.Item("y")
$wsY.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# Sheet "Q"
# ---------------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Range("C7").Value = 109.9450000000008
$wsQ.Range("C8").Value = 117.5900000000008
$wsQ.Range("C9").Value = 113.2700000000008
$wsQ.Range("C10").Value = 119.1550000000008
$wsQ.Range("C11").Value = 115.8050000000008
$wsQ.Range("C12").Value = 188.8550000000006
$wsQ.Range("C13").Value = 192.9200000000006
$wsQ.Range("C14").Value = 178.5050000000006
$wsQ.Range("C15").Value = 189.2700000000006
$wsQ.Range("C16").Value = 182.1250000000006
$wsQ.Range("C17").Value = 40.35
$wsQ.Range("C18").Value = 30.90499999999942
$wsQ.Range("C19").Value = 27.59499999999942
$wsQ.Range("C20").Value = 31.97499999999942
$wsQ.Range("C21").Value = 33.99499999999941
$wsQ.Range("C22").Value = 112.9799999999989
$wsQ.Range("C23").Value = 110.9599999999989
$wsQ.Range("C24").Value = 109.9149999999989
$wsQ.Range("C25").Value = 115.9299999999989
$wsQ.Range("C26").Value = 114.809999999999
$wsQ.Range("C27").Value = 258.7350000000008
$wsQ.Range("C28").Value = 269.2400000000008
$wsQ.Range("C29").Value = 250.9150000000008
$wsQ.Range("C30").Value = 261.9150000000008
$wsQ.Range("C31").Value = 255.0150000000008
$wsQ.Range("C32").Value = 107.3799999999999
$wsQ.Range("C33").Value = 112.2399999999999
$wsQ.Range("C34").Value = 93.78999999999985
$wsQ.Range("C35").Value = 108.8349999999998
$wsQ.Range("C36").Value = 94.77999999999986
$wsQ.Range("C37").Value = 193.0200000000017
$wsQ.Range("C38").Value = 202.3100000000017
$wsQ.Range("C39").Value = 191.2450000000017
$wsQ.Range("C40").Value = 208.9250000000017
$wsQ.Range("C41").Value = 197.6600000000017
$wsQ.Range("C42").Value = 180.9049999999984
$wsQ.Range("C43").Value = 190.1199999999983
$wsQ.Range("C44").Value = 169.7349999999983
$wsQ.Range("C45").Value = 179.7
$wsQ.Range("C46").Value = 173.7399999999984
$wsQ.Range("C47").Value = 266.3899999999988
$wsQ.Range("C48").Value = 278.0849999999988
$wsQ.Range("C49").Value = 249.45
$wsQ.Range("C50").Value = 270.4299999999989
$wsQ.Range("C51").Value = 258.4699999999989
$wsQ.Range("C52").Value = 250.970000000001
$wsQ.Range("C53").Value = 260.9900000000009
$wsQ.Range("C54").Value = 252.975000000001
$wsQ.Range("C55").Value = 269.580000000001
$wsQ.Range("C56").Value = 250.575000000001
$wsQ.Range("C57").Value = 250.970000000001
$wsQ.Range("C58").Value = 260.9900000000009
$wsQ.Range("C59").Value = 252.975000000001
$wsQ.Range("C60").Value = 269.580000000001
$wsQ.Range("C61").Value = 250.575000000001
$wsQ.Range("C62").Value = 258.7350000000008
$wsQ.Range("C63").Value = 269.2400000000008
$wsQ.Range("C64").Value = 250.9150000000008
$wsQ.Range("C65").Value = 261.9150000000008
$wsQ.Range("C66").Value = 255.0150000000008
$wsQ.Range("C67").Value = 266.3899999999988
$wsQ.Range("C68").Value = 278.0849999999988
$wsQ.Range("C69").Value = 249.45
$wsQ.Range("C70").Value = 270.4299999999989
$wsQ.Range("C71").Value = 258.4699999999989

# ---------------------------------------------------------------------------
# Sheet "R"
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("R")
$wsR.Range("C3").Value = 0
$wsR.Range("C5").Value = 0
$wsR.Range("C12").Value = 6.57
$wsR.Range("C13").Value = 5.2
$wsR.Range("C14").Value = 7.32
$wsR.Range("C15").Value = 5.515
$wsR.Range("C16").Value = 5.44

# ---------------------------------------------------------------------------
# Sheet "L"
# ---------------------------------------------------------------------------
$wsL = $wb.Worksheets.Item("L")
$wsL.Range("C22").Value = 0
$wsL.Range("C23").Value = 0
$wsL.Range("C24").Value = 0
$wsL.Range("C25").Value = 0
$wsL.Range("C26").Value = 0
$wsL.Range("C42").Value = 0
$wsL.Range("C43").Value = 0
$wsL.Range("C44").Value = 0
$wsL.Range("C45").Value = 0
$wsL.Range("C46").Value = 0

# ---------------------------------------------------------------------------
# Sheet "rho": only the header row remains
# ---------------------------------------------------------------------------
$wsRho = $wb.Worksheets.Item("rho")
$wsRho.Rows("2:3").Delete()

# ---------------------------------------------------------------------------
# Sheet "alpha": only the header row remains
# ---------------------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Rows("2:3").Delete()
